$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (ifoCAST full series evaluation update)
$data = @{
    2  = @(-0.1025119037381019, 0.6613299526921388, 0.8215310669073134, 0.9063835098385856, 0.9266766292061702, 18)
    3  = @(0.007029451631606786, 0.5882322708542835, 0.4919477886271549, 0.7013898977224828, 0.7229398480521881, 17)
    4  = @(0.08761989786705919, 0.4400362220161121, 0.2943530511839517, 0.5425431330170457, 0.5529805828361719, 16)
    5  = @(0.126447536190323, 0.5254135465362988, 0.4807614626535618, 0.693369643591037, 0.7056702613561491, 15)
    6  = @(0.1753236501719707, 0.5133168387277518, 0.4122335504200531, 0.6420541647089076, 0.6409687592018587, 14)
    7  = @(0.1711200621906368, 0.5664235858473783, 0.4676572605208214, 0.6838547071716488, 0.6891345661333148, 13)
    8  = @(0.2215176019925603, 0.5017932260683228, 0.3910326139135445, 0.6253260061068502, 0.610778251015726, 12)
    9  = @(0.214449951378639, 0.4406311141385629, 0.2572337742857695, 0.5071821904264477, 0.4820471884612707, 11)
    10 = @(0.1725523049316589, 0.5360107823372657, 0.4134394327613503, 0.6429925604245746, 0.6529123940078213, 10)
    11 = @(0.08702659289064565, 0.552660861616738, 0.4356818832333129, 0.660062029837585, 0.693989760215202, 9)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
